$wb = $excel.ActiveWorkbook

# --- Add the new "mapping" sheet after "target" (do this first so the new
#     shared strings "source column"/"target column" get lower sst indices
#     than the "*_target" field names, matching the canonical ordering) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$mapping = $wb.Worksheets.Add($null, $lastSheet)
$mapping.Name = "mapping"

# Column widths (closest values reachable through the host's width-rounding grid)
$mapping.Columns.Item(2).ColumnWidth = 13.8
$mapping.Columns.Item(3).ColumnWidth = 20.5

# Header row
$mapping.Range("A1").Value = "file_name"
$mapping.Range("B1").Value = "source column"
$mapping.Range("C1").Value = "target column"

# --- Update the "target" sheet's column B values (field -> field_target) ---
$target = $wb.Worksheets.Item("target")
$target.Range("B2").Value = "emp_id_target"
$target.Range("B3").Value = "emp_name_target"
$target.Range("B4").Value = "order_id_target"
$target.Range("B5").Value = "order_status_target"

# Move the selection on the target sheet (it is no longer the active/selected tab)
[void]$target.Range("B8").Select()

# Data rows: file name, source field (via formula to source sheet), target field (via formula to target sheet)
$mapping.Range("A2").Value = "file1"
$mapping.Range("B2").Formula = "=source!B2"
$mapping.Range("C2").Formula = "=target!B2"

$mapping.Range("A3").Value = "file1"
$mapping.Range("B3").Formula = "=source!B3"
$mapping.Range("C3").Formula = "=target!B3"

$mapping.Range("A4").Value = "file2"
$mapping.Range("B4").Formula = "=source!B4"
$mapping.Range("C4").Formula = "=target!B4"

$mapping.Range("A5").Value = "file2"
$mapping.Range("B5").Formula = "=source!B5"
$mapping.Range("C5").Formula = "=target!B5"

[void]$mapping.Range("C6").Select()
